$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 182 ---
$ws.Range('A182').Value = '2026-02-15 05:51:59'
$ws.Range('C182').Value = 'N A DAMA ADAMA'
$ws.Range('D182').Value = 354745

# --- Append new rows 231-265 ---
$ws.Range("B231:B265").NumberFormat = "@"

$ws.Range('A231').Value = '2026-02-15 05:49:21'
$ws.Range('B231').Value = '237671646117'
$ws.Range('C231').Value = 'FLORENCE MEDONGUE'
$ws.Range('D231').Value = 623824

$ws.Range('A232').Value = '2026-02-15 05:12:38'
$ws.Range('B232').Value = '237674440808'
$ws.Range('C232').Value = 'HOUMI EPSE MVEINGUE GUINDOP TATIANA ETS LE CONTENT'
$ws.Range('D232').Value = 6322

$ws.Range('A233').Value = '2026-02-15 05:53:56'
$ws.Range('B233').Value = '237675788721'
$ws.Range('C233').Value = 'ETS MOBILE FINANCIAL SERVICES MFS TCHAWE MBOUGA JUDITH FLORE'
$ws.Range('D233').Value = 569488

$ws.Range('A234').Value = '2026-02-15 05:19:20'
$ws.Range('B234').Value = '237676036914'
$ws.Range('C234').Value = 'NDENGUE ELOKO SAMUEL HERVE ETS MOBILE FINANCIAL SERVICES MFS'
$ws.Range('D234').Value = 92721

$ws.Range('A235').Value = '2026-02-14 13:31:12'
$ws.Range('B235').Value = '237679550294'
$ws.Range('C235').Value = 'N A FONATIA'
$ws.Range('D235').Value = 35722

$ws.Range('A236').Value = '2026-02-14 10:55:51'
$ws.Range('B236').Value = '237679604574'
$ws.Range('C236').Value = 'KINGUE KOMBI VICTORINE SIDONIE VISION TRADING COMPLEX AND TECHNOLOGIES SARL VISION TRADING COMPLEX'
$ws.Range('D236').Value = 25643

$ws.Range('A237').Value = '2026-02-14 14:38:23'
$ws.Range('B237').Value = '237682370358'
$ws.Range('C237').Value = 'JEANNE PRISCA NGO DJON EPSE EBANA ZOE'
$ws.Range('D237').Value = 328505

$ws.Range('A238').Value = '2026-02-15 05:56:21'
$ws.Range('B238').Value = '237682639044'
$ws.Range('C238').Value = 'JOSEPH KAMGA'
$ws.Range('D238').Value = 46373

$ws.Range('A239').Value = '2026-02-15 04:58:37'
$ws.Range('B239').Value = '237683023087'
$ws.Range('C239').Value = 'FAGHUIE ABIBA'
$ws.Range('D239').Value = 373215

$ws.Range('A240').Value = '2026-02-15 05:49:11'
$ws.Range('B240').Value = '237650874464'
$ws.Range('C240').Value = 'ASSONFACK VANESSA ATB POINT COM'
$ws.Range('D240').Value = 536391

$ws.Range('A241').Value = '2026-02-15 05:54:56'
$ws.Range('B241').Value = '237653854849'
$ws.Range('C241').Value = 'TITTI GASTON CLEMENT TOP MOBIL'
$ws.Range('D241').Value = 98719

$ws.Range('A242').Value = '2026-02-14 11:45:07'
$ws.Range('B242').Value = '237654164073'
$ws.Range('C242').Value = 'NGUIDJOL SIMONE ASTRIDE SPECTRUM SPECTRUM'
$ws.Range('D242').Value = 54178

$ws.Range('A243').Value = '2026-02-14 16:47:25'
$ws.Range('B243').Value = '237671357520'
$ws.Range('C243').Value = 'CALICE WOTI EPSE DJOMO'
$ws.Range('D243').Value = 15906

$ws.Range('A244').Value = '2026-02-13 13:41:17'
$ws.Range('B244').Value = '237672587687'
$ws.Range('C244').Value = 'Kamaha Tomy Nadine LA NEGRESSE SARL'
$ws.Range('D244').Value = 247714

$ws.Range('A245').Value = '2026-02-15 05:50:59'
$ws.Range('B245').Value = '237674240552'
$ws.Range('C245').Value = 'AUGUSTINE NGO BAYOI'
$ws.Range('D245').Value = 56178

$ws.Range('A246').Value = '2026-02-15 03:52:03'
$ws.Range('B246').Value = '237675239360'
$ws.Range('C246').Value = 'ERIC MBAH AKEN'
$ws.Range('D246').Value = 107033

$ws.Range('A247').Value = '2026-02-15 05:55:37'
$ws.Range('B247').Value = '237675396752'
$ws.Range('C247').Value = 'BENEDICTE CHANTAL MANTSANG'
$ws.Range('D247').Value = 121999

$ws.Range('A248').Value = '2026-02-15 05:23:39'
$ws.Range('B248').Value = '237675626141'
$ws.Range('C248').Value = 'FLORENCE NGUEFACK'
$ws.Range('D248').Value = 148545

$ws.Range('A249').Value = '2026-02-15 02:12:28'
$ws.Range('B249').Value = '237676840777'
$ws.Range('C249').Value = 'ETP109 ETP'
$ws.Range('D249').Value = 0

$ws.Range('A250').Value = '2026-02-15 05:43:40'
$ws.Range('B250').Value = '237677833877'
$ws.Range('C250').Value = 'ISSA ISSYAKOU'
$ws.Range('D250').Value = 561227

$ws.Range('A251').Value = '2026-02-15 05:25:47'
$ws.Range('B251').Value = '237678854978'
$ws.Range('C251').Value = 'NSAMO NDJOUOHOU MICRANGE ETS MOBILE FINANCIAL SERVICES MFS'
$ws.Range('D251').Value = 258563

$ws.Range('A252').Value = '2026-02-15 05:00:13'
$ws.Range('B252').Value = '237679422591'
$ws.Range('C252').Value = 'ETS LE CONTENT 42'
$ws.Range('D252').Value = 401400

$ws.Range('A253').Value = '2026-02-15 06:04:06'
$ws.Range('B253').Value = '237650353920'
$ws.Range('C253').Value = 'MENIAPI HELENE EDOSSINE TOP MOBIL TELECOM'
$ws.Range('D253').Value = 980204

$ws.Range('A254').Value = '2026-02-14 13:20:06'
$ws.Range('B254').Value = '237651927448'
$ws.Range('C254').Value = 'charity aben awalah'
$ws.Range('D254').Value = 71305

$ws.Range('A255').Value = '2026-02-14 12:20:24'
$ws.Range('B255').Value = '237653294562'
$ws.Range('C255').Value = 'NANHOU KEMAYOU AVIGAEL ETS MOBILE FINANCIAL SERVICES MFS'
$ws.Range('D255').Value = 271672

$ws.Range('A256').Value = '2026-02-14 13:53:49'
$ws.Range('B256').Value = '237678046498'
$ws.Range('C256').Value = 'MFS SOCAVER'
$ws.Range('D256').Value = 303

$ws.Range('A257').Value = '2026-02-14 14:41:18'
$ws.Range('B257').Value = '237679428698'
$ws.Range('C257').Value = 'ETS LE CONTENT 29'
$ws.Range('D257').Value = 7

$ws.Range('A258').Value = '2026-02-14 17:33:53'
$ws.Range('B258').Value = '237679551262'
$ws.Range('C258').Value = 'LA NEGRESSE LTDLA CBOX R1 MEGNE JUDITH'
$ws.Range('D258').Value = 17349

$ws.Range('A259').Value = '2026-02-14 11:43:01'
$ws.Range('B259').Value = '237680574202'
$ws.Range('C259').Value = 'TOUMEWO SAMUEL'
$ws.Range('D259').Value = 297874

$ws.Range('A260').Value = '2026-02-15 05:13:29'
$ws.Range('B260').Value = '237681118330'
$ws.Range('C260').Value = 'SAHA NDESA JONAS LTDLA_POLAS_OTH_NDOGBONG SERIE'
$ws.Range('D260').Value = 303554

$ws.Range('A261').Value = '2026-02-15 06:04:31'
$ws.Range('B261').Value = '237674446293'
$ws.Range('C261').Value = 'SYDONIE MAFOMA MESSINE'
$ws.Range('D261').Value = 10536

$ws.Range('A262').Value = '2026-02-15 04:29:58'
$ws.Range('B262').Value = '237679085953'
$ws.Range('C262').Value = 'MADELEINE NKOUADJIO'
$ws.Range('D262').Value = 24018

$ws.Range('A263').Value = '2026-02-15 05:47:17'
$ws.Range('B263').Value = '237681662761'
$ws.Range('C263').Value = 'EMMANUEL EKOLLE ELUMBA'
$ws.Range('D263').Value = 32326

$ws.Range('A264').Value = '2026-02-15 04:43:22'
$ws.Range('B264').Value = '237682975726'
$ws.Range('C264').Value = 'SYLVIE-ISABELLE DGANHOU EPSE KOUAHOU'
$ws.Range('D264').Value = 76155

$ws.Range('A265').Value = '2026-02-15 01:45:11'
$ws.Range('B265').Value = '237683075075'
$ws.Range('C265').Value = 'ESSOM YOUASSI FRANCK LIONEL STYLE. COM'
$ws.Range('D265').Value = 4943

$ws.Range("B231:B265").Style = "Normal"
